$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.5022738159769057
$ws.Range("J2").Value = 0.5022738159769057
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 1.084983566359333
$ws.Range("R2").Value = 9.764852097234
$ws.Range("S2").Value = 0.01057599242529592
$ws.Range("T2").Value = 0.01057599242529592

# Row 3
$ws.Range("I3").Value = 0.5022738159769057
$ws.Range("J3").Value = 0.5022738159769057
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("S3").Value = 0.3884069259558694
$ws.Range("T3").Value = 0.3884069259558695

# Row 4
$ws.Range("I4").Value = 0.5022738159769057
$ws.Range("J4").Value = 0.5022738159769057
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 0.276592355238
$ws.Range("R4").Value = 2.489331197142
$ws.Range("S4").Value = 0.002696113327971865
$ws.Range("T4").Value = 0.002696113327971866

# Row 5
$ws.Range("I5").Value = 0.5022738159769057
$ws.Range("J5").Value = 0.5022738159769057
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 10.182426744768
$ws.Range("R5").Value = 91.64184070291201
$ws.Range("S5").Value = 0.09925428500742062
$ws.Range("T5").Value = 0.09925428500742063

# Row 6
$ws.Range("I6").Value = 0.5022738159769057
$ws.Range("J6").Value = 0.5022738159769057
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 0.1375208689366667
$ws.Range("R6").Value = 1.23768782043
$ws.Range("S6").Value = 0.001340499260347887
$ws.Range("T6").Value = 0.001340499260347887

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.63173
$ws.Range("H7").Value = 1.89519
$ws.Range("I7").Value = 0.4977261840230943
$ws.Range("J7").Value = 0.4977261840230943
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 1.075160028323333
$ws.Range("R7").Value = 9.67644025491
$ws.Range("S7").Value = 0.0104802364460538
$ws.Range("T7").Value = 0.0104802364460538

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.63173
$ws.Range("H8").Value = 1.89519
$ws.Range("I8").Value = 0.4977261840230943
$ws.Range("J8").Value = 0.4977261840230943
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("Q8").Value = 39.48561843831001
$ws.Range("R8").Value = 355.37056594479
$ws.Range("S8").Value = 0.3848902549860257
$ws.Range("T8").Value = 0.3848902549860258

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.63173
$ws.Range("H9").Value = 1.89519
$ws.Range("I9").Value = 0.4977261840230943
$ws.Range("J9").Value = 0.4977261840230943
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 0.27408806337
$ws.Range("R9").Value = 2.46679257033
$ws.Range("S9").Value = 0.002671702477293666
$ws.Range("T9").Value = 0.002671702477293667

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.63173
$ws.Range("H10").Value = 1.89519
$ws.Range("I10").Value = 0.4977261840230943
$ws.Range("J10").Value = 0.4977261840230943
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 10.09023414432
$ws.Range("R10").Value = 90.81210729888001
$ws.Range("S10").Value = 0.09835562785330529
$ws.Range("T10").Value = 0.0983556278533053

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.63173
$ws.Range("H11").Value = 1.89519
$ws.Range("I11").Value = 0.4977261840230943
$ws.Range("J11").Value = 0.4977261840230943
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 0.1362757427166667
$ws.Range("R11").Value = 1.22648168445
$ws.Range("S11").Value = 0.001328362260415765
$ws.Range("T11").Value = 0.001328362260415765

